$d = $word.ActiveDocument

# =========================================================================
# Paragraph: "...Bovendien hebben we een breadth solver en een random
#             solver..." paragraph.
# =========================================================================

# 1) Merge "spelen." onto the end of the preceding run, so it reads
#    "...welke game je wilt spelen. " (instead of being its own run).
$rngInsert = $d.Content
$rngInsert.Find.Execute("kiezen welke game je wilt ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rngInsert.Collapse(0)
$rngInsert.InsertAfter("spelen. ")

# 2) Remove the now-redundant original "spelen." run (search only after the
#    text we just inserted, so we target the old run and not the new text).
$rngOldWord = $d.Range($rngInsert.End, $d.Content.End)
$rngOldWord.Find.Execute("spelen.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rngOldWord.Text = ""

# 3) Remove the leftover single space that used to separate "spelen." from
#    "Bovendien hebben we een ".
$rngGap = $d.Content
$rngGap.Find.Execute(" Bovendien hebben we een", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rngSpace = $d.Range($rngGap.Start, $rngGap.Start + 1)
$rngSpace.Text = ""

# 4) "bij de random solver wordt de oplossing per solve" ->
#    "bij de random solver hebben we het zo geïmplementeerd dat de oplossing per solve"
$d.Content.Find.Execute(" wordt de oplossing per ", $true, $false, $false, $false, $false, $true, 1, $false, " hebben we het zo geïmplementeerd dat de oplossing per ", 2)

# 5) "solve steeds beter, ook worden er geen repetitions" ->
#    "solve steeds beter wordt en dat er geen repetitions"
$d.Content.Find.Execute(" steeds beter, ook worden er geen ", $true, $false, $false, $false, $false, $true, 1, $false, " steeds beter wordt en dat er geen ", 2)

# 6) "repetitions uitgevoerd." -> "repetitions worden uitgevoerd."
$d.Content.Find.Execute(" uitgevoerd. ", $true, $false, $false, $false, $false, $true, 1, $false, " worden uitgevoerd. ", 2)

# =========================================================================
# Paragraph: "We willen graag verder met de breadth solver, ..." ->
#            "We willen graag verder vanuit het idee van een breadth
#             solver, ..."
# =========================================================================
$d.Content.Find.Execute("We willen graag verder met de ", $true, $false, $false, $false, $false, $true, 1, $false, "We willen graag verder vanuit het idee van een ", 2)
